# lecture-2-practice-of-conversion-string-number-and-datetime:
# remove the content from datetime sheet
#
# Slide 18 ("All in datetime") body placeholder lists several topics as
# separate paragraphs. Remove the "Benchmarking" paragraph entirely,
# leaving the other paragraphs (and their formatting) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

for ($i = $tr.Paragraphs().Count; $i -ge 1; $i--) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.Trim() -eq "Benchmarking") {
        $para.Delete()
    }
}
